$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newQuery = @'
SELECT
    COUNT(DISTINCT std.dbgap_accession) AS "Studies",
    COUNT(DISTINCT prt.participant_id) AS "Participants",
    COUNT(DISTINCT smp.sample_id) AS "Samples",
    (COUNT(DISTINCT seq.sequencing_file_id) + COUNT(DISTINCT maf.methylation_array_file_id)) AS "Files"
FROM 
    df_study std
LEFT JOIN df_participant prt ON std.id = prt."study.id"
LEFT JOIN df_sample smp ON prt.id = smp."participant.id"
LEFT JOIN df_diagnosis dgn ON prt.id = dgn."participant.id"
LEFT JOIN df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN df_sequencing_file seq ON smp.id = seq."sample.id"
LEFT JOIN df_treatment  trt ON prt.id = trt."participant_id"
LEFT JOIN df_methylation_array_file maf ON smp.id = maf."sample.id"
WHERE 
std.dbgap_accession = 'phs002504' 
AND prt.sex_at_birth = 'Female'
AND prt.race = 'Unknown'
AND trt.treatment_type LIKE '%Chemotherapy%';
'@

# Here-strings retain a trailing newline from the closing delimiter line; strip it
$newQuery = $newQuery.TrimEnd("`r", "`n")

$ws.Range("C2").Value = $newQuery
